# feat: add 2022-Q3 data
#
# The workbook tracks quarterly fund-holding snapshots for a stock. A new
# quarter ("2022-Q3") is published, so:
#   1. A new sheet is inserted right after "总计" (pushing the older
#      quarters back) by duplicating the most-recent quarter sheet
#      ("2022-Q2") and overwriting it with the new quarter's figures.
#   2. The "总计" (totals) roll-up sheet gets a new leading row for
#      2022-Q3 and all the other rows shift down by one.

$wb = $excel.ActiveWorkbook

# --- 1. Add the new "2022-Q3" sheet right after "总计" -------------------
$totalSheet = $wb.Worksheets.Item("总计")
$prevQuarterSheet = $wb.Worksheets.Item("2022-Q2")

# Duplicate the latest quarter's sheet (keeps headers/styles/column widths)
# and drop the copy right after "总计"; the original "2022-Q2" sheet is
# left untouched further down the tab order.
$prevQuarterSheet.Copy($null, $totalSheet)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q3"

# Overwrite the copied row with the new quarter's fund data.
# (B2 fund code "010764" is unchanged from the copied source, so it's left
# alone — re-assigning it as a plain string would strip the leading zero.)
$newSheet.Range("C2").Value = "九泰锐升混合"

$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value = "2.15"

$newSheet.Range("E2").NumberFormat = "@"
$newSheet.Range("E2").Value = "54.42"

$newSheet.Range("F2").NumberFormat = "@"
$newSheet.Range("F2").Value = "2.56"

$newSheet.Range("G2").NumberFormat = "@"
$newSheet.Range("G2").Value = "0.0550"

$newSheet.Range("H2").Value = 8

# --- 2. Insert a new leading row in "总计" for 2022-Q3 --------------------
# Shift the existing quarter rows (2-5) down to rows 3-6, values + styles
# together, then restore row 2's index-column style (the shift leaves it
# blank) before filling in the new quarter and renumbering the index col.
$totalSheet.Range("A2:D5").Copy($totalSheet.Range("A3:D6"))

$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0.06
